$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "tentativa de escrever no pach vermelho" - change A2's text and paint A1 red
$ws.Range("A2").Value = "vai"
$ws.Range("A1").Interior.Color = 255
